# Insert a new data row at row 90 (pushing existing rows 90-104 down to 91-105)
# and populate it with the new "Poroto verde" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(90).Insert()

$ws.Range("A90").Value = 5
$ws.Range("B90").Value = 'Macroferia Regional de Talca'
$ws.Range("C90").Value = 'Maule'
$ws.Range("D90").Value = 44505
$ws.Range("E90").Value = 7
$ws.Range("F90").Value = 100112031
$ws.Range("G90").Value = 'Poroto verde'
$ws.Range("H90").Value = 'Sin especificar'
$ws.Range("I90").Value = 'Primera'
$ws.Range("J90").Value = 150
$ws.Range("K90").Value = 30000
$ws.Range("L90").Value = 30000
$ws.Range("M90").Value = 30000
$ws.Range("N90").Value = '$/saco 25 kilos'
$ws.Range("O90").Value = 'Región del Maule'
$ws.Range("P90").Value = 1200
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = 'Hortaliza'
